$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual cell values as described by the diff
$ws.Range("E5").Value  = 13.58789999999999
$ws.Range("E6").Value  = 11.685
$ws.Range("D7").Value  = -7.837099999999992
$ws.Range("E7").Value  = 12.1008
$ws.Range("A8").Value  = -21.11560000000002
$ws.Range("E8").Value  = 12.51320000000001
$ws.Range("E9").Value  = 9.650599999999988
$ws.Range("A10").Value = -20.48599999999997
$ws.Range("E10").Value = 11.28369999999999
$ws.Range("A12").Value = -22.44490000000003
$ws.Range("E12").Value = 12.76569999999999
$ws.Range("B13").Value = 6.343599999999998
$ws.Range("A18").Value = -22.38120000000002
$ws.Range("D20").Value = -8.278699999999999
$ws.Range("A25").Value = -22.35920000000004
